$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0374324806034565
$ws.Range("C2").Value = 0.010669447481632233
$ws.Range("D2").Value = 0.006921728607267141
$ws.Range("E2").Value = 0.004782211035490036
$ws.Range("F2").Value = 0.0000001204314088454339
$ws.Range("G2").Value = 0.00132271577604115
$ws.Range("J2").Value = 0.12753915786743164
$ws.Range("K2").Value = 1.4484678506851196
